$wb = $excel.ActiveWorkbook

# Access sheets by position to avoid ambiguity between case-variant names
# (e.g. "Vector_bf" vs "Vector_BF" resolve the same way when looked up by name).
# 1 Funciones_Objetivo
# 2 Restricciones_del_lider
# 3 Restricciones_del_follower
# 4 Punto_modificado
# 5 Vector_bf
# 6 Vector_BF
# 7 Vector_Alpha
$ws2 = $wb.Worksheets.Item(2)
$ws3 = $wb.Worksheets.Item(3)
$ws4 = $wb.Worksheets.Item(4)
$ws5 = $wb.Worksheets.Item(5)
$ws6 = $wb.Worksheets.Item(6)
$ws7 = $wb.Worksheets.Item(7)

# Ranges that will receive text-like numeric strings. Force them to Text
# number format first so Excel stores them as strings (matching the
# original workbook's representation) rather than converting to numbers.
$r2 = $ws2.Range("A2:D3")
$r3 = $ws3.Range("A2:F3")
$r4 = $ws4.Range("A2:B2")
$r5 = $ws5.Range("A2")
$r6 = $ws6.Range("A2:A3")

$r2.NumberFormat = "@"
$r3.NumberFormat = "@"
$r4.NumberFormat = "@"
$r5.NumberFormat = "@"
$r6.NumberFormat = "@"

# Restricciones_del_lider
$ws2.Range("A2").Value = "4.5 - x"
$ws2.Range("B2").Value = "-5.0"
$ws2.Range("D2").Value = "0.34"
$ws2.Range("A3").Value = "-4.5 + x"
$ws2.Range("B3").Value = "4.0"
$ws2.Range("D3").Value = "0.0"

# Restricciones_del_follower
$ws3.Range("A2").Value = "-3.880149812734083 + 1.3857677902621726y"
$ws3.Range("B2").Value = "2.880149812734083"
$ws3.Range("D2").Value = "0.09"
$ws3.Range("E2").Value = "9.5"
$ws3.Range("F2").Value = "3.7"
$ws3.Range("A3").Value = "4.675999999999999 - 1.67y"
$ws3.Range("B3").Value = "-5.675999999999999"
$ws3.Range("D3").Value = "0.82"
$ws3.Range("E3").Value = "3.5"
$ws3.Range("F3").Value = "0.4"

# Punto_modificado
$ws4.Range("A2").Value = "4.5"
$ws4.Range("B2").Value = "2.8"

# Vector_bf
$ws5.Range("A2").Value = "-3.255319101123596"

# Vector_BF
$ws6.Range("A2").Value = "-1.3299999999999998"
$ws6.Range("A3").Value = "-6.31979400749064"

# Restore the default (General) number format now that the text values are set,
# so the cells keep the same style/formatting as the rest of the workbook.
$r2.ClearFormats()
$r3.ClearFormats()
$r4.ClearFormats()
$r5.ClearFormats()
$r6.ClearFormats()

# Vector_Alpha (numeric cell, not a shared string)
$ws7.Range("A2").Value = 2.67
